$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure changed cells keep their original text formatting (e.g. trailing
# zeros, dot-grouped numbers, percent strings) instead of being auto-converted
# to numeric values by Excel when the new value looks like a number.
$cellsToUpdate = @(
    @{ Cell = 'D2'; Value = '61.876.05' }
    @{ Cell = 'E2'; Value = '  -1.08%  ' }
    @{ Cell = 'D3'; Value = '3.394.71' }
    @{ Cell = 'E3'; Value = '  -2.12%  ' }
    @{ Cell = 'E4'; Value = '  -0.19%  ' }
    @{ Cell = 'D5'; Value = '403.62' }
    @{ Cell = 'E5'; Value = '  -2.35%  ' }
    @{ Cell = 'D6'; Value = '132.18' }
    @{ Cell = 'E6'; Value = '  +2.73%  ' }
    @{ Cell = 'D7'; Value = '0.588' }
    @{ Cell = 'E7'; Value = '  -0.50%  ' }
    @{ Cell = 'D8'; Value = '1.00' }
    @{ Cell = 'E8'; Value = '  -0.13%  ' }
    @{ Cell = 'E9'; Value = '  -3.37%  ' }
    @{ Cell = 'E10'; Value = '  -5.19%  ' }
    @{ Cell = 'D11'; Value = '42.30' }
    @{ Cell = 'E11'; Value = '  -0.05%  ' }
    @{ Cell = 'E12'; Value = '  -1.78%  ' }
    @{ Cell = 'D13'; Value = '3.896.34' }
    @{ Cell = 'E13'; Value = '  -3.68%  ' }
    @{ Cell = 'D14'; Value = '8.39' }
    @{ Cell = 'E14'; Value = '  -3.50%  ' }
    @{ Cell = 'E15'; Value = '  -1.04%  ' }
    @{ Cell = 'D16'; Value = '3.378.88' }
    @{ Cell = 'E16'; Value = '  -3.35%  ' }
    @{ Cell = 'D17'; Value = '61.798.89' }
    @{ Cell = 'E17'; Value = '  -1.38%  ' }
    @{ Cell = 'E18'; Value = '  -2.20%  ' }
    @{ Cell = 'E19'; Value = '  +1.12%  ' }
    @{ Cell = 'E20'; Value = '  -4.99%  ' }
    @{ Cell = 'E21'; Value = '  -4.43%  ' }
    @{ Cell = 'D22'; Value = '84.40' }
    @{ Cell = 'E22'; Value = '  +2.95%  ' }
    @{ Cell = 'D23'; Value = '316.38' }
    @{ Cell = 'E23'; Value = '  +0.73%  ' }
    @{ Cell = 'D24'; Value = '12.71' }
    @{ Cell = 'E25'; Value = '  -2.46%  ' }
    @{ Cell = 'E26'; Value = '  +9.44%  ' }
    @{ Cell = 'D27'; Value = '29.54' }
    @{ Cell = 'E27'; Value = '  -3.64%  ' }
    @{ Cell = 'D28'; Value = '8.20' }
    @{ Cell = 'E28'; Value = '  +1.98%  ' }
    @{ Cell = 'D29'; Value = '7.73' }
    @{ Cell = 'E29'; Value = '  +0.39%  ' }
    @{ Cell = 'D30'; Value = '2.71' }
    @{ Cell = 'E30'; Value = '  +2.35%  ' }
    @{ Cell = 'E31'; Value = '  -3.93%  ' }
    @{ Cell = 'E32'; Value = '  -2.18%  ' }
    @{ Cell = 'B33'; Value = 'InjectiveProtocol' }
    @{ Cell = 'C33'; Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj' }
    @{ Cell = 'D33'; Value = '41.83' }
    @{ Cell = 'E33'; Value = '  -2.02%  ' }
    @{ Cell = 'B34'; Value = 'Dai' }
    @{ Cell = 'C34'; Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai' }
    @{ Cell = 'D34'; Value = '1.00' }
    @{ Cell = 'E34'; Value = '  +0.04%  ' }
    @{ Cell = 'D35'; Value = '11.36' }
    @{ Cell = 'E35'; Value = '  -3.16%  ' }
    @{ Cell = 'D36'; Value = '0.0481' }
    @{ Cell = 'E36'; Value = '  -3.26%  ' }
    @{ Cell = 'D37'; Value = '51.79' }
    @{ Cell = 'E37'; Value = '  -0.82%  ' }
    @{ Cell = 'D38'; Value = '0.998' }
    @{ Cell = 'E38'; Value = '  -0.17%  ' }
    @{ Cell = 'D39'; Value = '3.43' }
    @{ Cell = 'E39'; Value = '  -2.92%  ' }
    @{ Cell = 'D40'; Value = '2.96' }
    @{ Cell = 'E40'; Value = '  -2.29%  ' }
    @{ Cell = 'D41'; Value = '139.03' }
    @{ Cell = 'E41'; Value = '  +1.95%  ' }
    @{ Cell = 'E42'; Value = '  -1.54%  ' }
    @{ Cell = 'E43'; Value = '  -1.17%  ' }
    @{ Cell = 'D44'; Value = '0.293' }
    @{ Cell = 'E44'; Value = '  +2.57%  ' }
    @{ Cell = 'E45'; Value = '  +0.67%  ' }
    @{ Cell = 'D46'; Value = '16.69' }
    @{ Cell = 'E46'; Value = '  -2.57%  ' }
    @{ Cell = 'E47'; Value = '  -0.75%  ' }
    @{ Cell = 'D48'; Value = '21.36' }
    @{ Cell = 'E48'; Value = '  -2.80%  ' }
    @{ Cell = 'D49'; Value = '2.121.94' }
    @{ Cell = 'E49'; Value = '  -4.66%  ' }
    @{ Cell = 'E50'; Value = '  -5.41%  ' }
    @{ Cell = 'D51'; Value = '1.86' }
    @{ Cell = 'E51'; Value = '  +1.34%  ' }
)

foreach ($item in $cellsToUpdate) {
    $range = $ws.Range($item.Cell)
    $range.NumberFormat = "@"
    $range.Value = $item.Value
}
